# Auto-generated PowerShell Word COM-interop script
# Inserts the "class / id 속성" study-note paragraphs after the
# "글자 정렬 위치" paragraph, matching the target OOXML diff exactly.

$d = $word.ActiveDocument

# Locate the anchor paragraph: the one whose text is "글자 정렬 위치"
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($ptext -eq "글자 정렬 위치`r") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) { throw "anchor paragraph not found" }

$count = 19

# Step 1: create $count empty paragraphs right after the anchor, preserving order.
$cursor = $d.Paragraphs($anchorIndex)
for ($i = 0; $i -lt $count; $i++) {
    $cursor.Range.InsertParagraphAfter()
    $cursor = $d.Paragraphs($anchorIndex + $i + 1)
}

# Step 2: fill each new paragraph with its exact OOXML content.
# InsertXML REPLACES the (now-empty, standalone) paragraph range with the
# payload, so every paragraph -- including ones that should stay blank --
# is written explicitly to avoid leftover inherited run formatting.
$target = $d.Paragraphs(($anchorIndex + 1))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@)

$target = $d.Paragraphs(($anchorIndex + 2))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>h</w:t>
      </w:r>
      <w:r>
        <w:t>tml</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">속성 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">class </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>속성</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 3))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>&lt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> class =”</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>aaa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>”&gt;123&lt;/a&gt;</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 4))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4050"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>&lt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> class =”</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>aaa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>bbb</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>”&gt;</w:t>
      </w:r>
      <w:r>
        <w:t>456</w:t>
      </w:r>
      <w:r>
        <w:t>&lt;/a&gt;</w:t>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>//</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>띄어쓰기로 여러 클래스 지정 가능</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 5))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>&lt;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> class </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>=”</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>bbb</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>”&gt;</w:t>
      </w:r>
      <w:r>
        <w:t>789</w:t>
      </w:r>
      <w:r>
        <w:t>&lt;/a&gt;</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 6))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">여러 개를 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>그루핑해서</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> 관리할 수 있다.</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 7))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:t>tyle</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">태그 </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>내부 에서</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> 호출은</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 8))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>aaa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 9))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 10))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@)

$target = $d.Paragraphs(($anchorIndex + 11))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>d 속성</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 12))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>하나의 개체에만 적용되며 중복은 불가능하다</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 13))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>&lt;</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>id</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> =”</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>aaa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>”&gt;</w:t>
      </w:r>
      <w:r>
        <w:t>000</w:t>
      </w:r>
      <w:r>
        <w:t>&lt;/a&gt;</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 14))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:t>tyle</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>태그 내부 호출은</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 15))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>#</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:t>aa{</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 16))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@)

$target = $d.Paragraphs(($anchorIndex + 17))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
'@)

$target = $d.Paragraphs(($anchorIndex + 18))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@)

$target = $d.Paragraphs(($anchorIndex + 19))
$target.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">적용 우선순위 </w:t>
      </w:r>
      <w:r>
        <w:t>id(#)&gt;class(.)&gt;</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>태그지정 이다</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@)

